$wb = $excel.ActiveWorkbook

# Rename the "quiz" sheet to "Quiz"
$quizSheet = $wb.Worksheets.Item("quiz")
$quizSheet.Name = "Quiz"

# Update the Training sheet header row (B1 becomes "Symptoms"; the rest of the
# header cells keep their existing text) and clear out the now-unused L column.
$ws = $wb.Worksheets.Item("Training")
$ws.Range("B1").Value = "Symptoms"
$ws.Range("L1:L78").ClearContents()

# Make Training the active sheet/tab, with B5 selected.
$ws.Activate()
$ws.Range("B5").Select()
